# Grade12StudentCredentials (sheet6 / rId6) gains a new auto-generated
# account row, matching the pattern already used on the other grade
# credential sheets (GradeOne/Grade4/Grade9StudentCredentials): a fresh
# "AutoYYYYMMDDHHMMSSmmm" username in column A and the shared default
# password "Password@123" in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade12StudentCredentials")

$ws.Range("A2").Value = "Auto20210827004426713"
$ws.Range("B2").Value = "Password@123"

# Give column A an explicit width like the sibling sheets (e.g.
# GradeOneStudentCredentials / Grade9StudentCredentials) do for their
# narrow "UserId" column.
$ws.Columns.Item(1).ColumnWidth = 6.3
